# The Evaluation Warning sheet (auto-generated by the Aspose.Cells library
# used to produce this fixture) stamps the current copyright year into a
# cell. Regenerating the fixture with a newer library build bumps that
# year from 2014 to 2016 - update the text in place so the existing
# shared-string slot / cell style are preserved.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Evaluation Warning")
$ws.Range("A5").Value = "Evaluation Only. Created with Aspose.Cells for .NET.Copyright 2003 - 2016 Aspose Pty Ltd."
